# Apply cryptos list update (prices/volumes refreshed, two coin rows swapped)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text formatting
# (values like "28.317.68" or "1.00" must stay literal text, not be
# reinterpreted as numbers/dates by Excel).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '28.317.68'
$ws.Range("E2").Value = '  -1.52%  '
$ws.Range("D3").Value = '1.552.42'
$ws.Range("E3").Value = '  -1.44%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '209.74'
$ws.Range("E6").Value = '  -1.93%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '23.73'
$ws.Range("E8").Value = '  -2.03%  '
$ws.Range("E9").Value = '  -2.12%  '
$ws.Range("D10").Value = '0.0584'
$ws.Range("D11").Value = '0.0890'
$ws.Range("E11").Value = '  +0.12%  '
$ws.Range("D12").Value = '1.774.43'
$ws.Range("D13").Value = '1.571.43'
$ws.Range("E13").Value = '  -0.11%  '
$ws.Range("D14").Value = '28.312.17'
$ws.Range("E14").Value = '  -1.46%  '
$ws.Range("E15").Value = '  -1.86%  '
$ws.Range("E16").Value = '  -2.49%  '
$ws.Range("D17").Value = '60.54'
$ws.Range("E17").Value = '  -3.13%  '
$ws.Range("D18").Value = '228.01'
$ws.Range("E18").Value = '  -1.07%  '
$ws.Range("E19").Value = '  -0.99%  '
$ws.Range("E20").Value = '  -2.93%  '
$ws.Range("D21").Value = '1.00'
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("E22").Value = '  +0.52%  '
$ws.Range("D23").Value = '8.89'
$ws.Range("E23").Value = '  -3.21%  '
$ws.Range("E24").Value = '  -4.28%  '
$ws.Range("D25").Value = '151.33'
$ws.Range("E25").Value = '  -0.36%  '
$ws.Range("D26").Value = '14.75'
$ws.Range("E26").Value = '  -1.99%  '
$ws.Range("E28").Value = '  -0.08%  '
$ws.Range("E29").Value = '  -3.30%  '
$ws.Range("D30").Value = '0.0468'
$ws.Range("E30").Value = '  -3.66%  '
$ws.Range("E31").Value = '  -4.61%  '
$ws.Range("E32").Value = '  -1.57%  '
$ws.Range("B33").Value = 'Maker'
$ws.Range("C33").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D33").Value = '1.385.88'
$ws.Range("E33").Value = '  -0.29%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").Value = '3.03'
$ws.Range("E34").Value = '  -2.98%  '
$ws.Range("E35").Value = '  +1.00%  '
$ws.Range("E36").Value = '  -3.67%  '
$ws.Range("D37").Value = '2.34'
$ws.Range("E37").Value = '  -1.23%  '
$ws.Range("E38").Value = '  -0.99%  '
$ws.Range("E39").Value = '  -3.20%  '
$ws.Range("E40").Value = '  +1.14%  '
$ws.Range("D41").Value = '0.512'
$ws.Range("E41").Value = '  -2.93%  '
$ws.Range("E42").Value = '  -0.10%  '
$ws.Range("E43").Value = '  -2.35%  '
$ws.Range("D44").Value = '0.0467'
$ws.Range("E44").Value = '  -0.70%  '
$ws.Range("E45").Value = '  -2.50%  '
$ws.Range("D46").Value = '61.91'
$ws.Range("E46").Value = '  -2.38%  '
$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D47").Value = '1.687.85'
$ws.Range("E47").Value = '  -1.40%  '
$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").Value = '0.906'
$ws.Range("E48").Value = '  -6.36%  '
$ws.Range("D49").Value = '85.73'
$ws.Range("E49").Value = '  -1.16%  '
$ws.Range("D50").Value = '42.06'
$ws.Range("E50").Value = '  +5.54%  '
$ws.Range("E51").Value = '  +0.58%  '
